$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.424164
$ws.Range("H2").Value = 4.272492
$ws.Range("I2").Value = 0.03823384556371837
$ws.Range("J2").Value = 0.03823384556371837
$ws.Range("M2").Value = 19.42991633333333
$ws.Range("N2").Value = 58.289749
$ws.Range("O2").Value = 0.08673502554925175
$ws.Range("P2").Value = 0.08673502554925173
$ws.Range("Q2").Value = 27.67138736494534
$ws.Range("R2").Value = 249.042486284508
$ws.Range("S2").Value = 0.003316213571815259
$ws.Range("T2").Value = 0.003316213571815258
$ws.Range("G3").Value = 1.424164
$ws.Range("H3").Value = 4.272492
$ws.Range("I3").Value = 0.03823384556371837
$ws.Range("J3").Value = 0.03823384556371837
$ws.Range("O3").Value = 0.04103322570207864
$ws.Range("P3").Value = 0.04103322570207864
$ws.Range("Q3").Value = 13.09097767649467
$ws.Range("R3").Value = 117.818799088452
$ws.Range("S3").Value = 0.001568858014474474
$ws.Range("T3").Value = 0.001568858014474474
$ws.Range("G4").Value = 1.424164
$ws.Range("H4").Value = 4.272492
$ws.Range("I4").Value = 0.03823384556371837
$ws.Range("J4").Value = 0.03823384556371837
$ws.Range("M4").Value = 105.042315
$ws.Range("N4").Value = 315.126945
$ws.Range("O4").Value = 0.4689082402093144
$ws.Range("P4").Value = 0.4689082402093144
$ws.Range("Q4").Value = 149.59748349966
$ws.Range("R4").Value = 1346.37735149694
$ws.Range("S4").Value = 0.01792816523971789
$ws.Range("T4").Value = 0.01792816523971788
$ws.Range("G5").Value = 1.424164
$ws.Range("H5").Value = 4.272492
$ws.Range("I5").Value = 0.03823384556371837
$ws.Range("J5").Value = 0.03823384556371837
$ws.Range("M5").Value = 2.834125333333334
$ws.Range("N5").Value = 8.502376000000002
$ws.Range("O5").Value = 0.01265151784388958
$ws.Range("P5").Value = 0.01265151784388958
$ws.Range("Q5").Value = 4.036259271221335
$ws.Range("R5").Value = 36.326333440992
$ws.Range("S5").Value = 0.0004837161793899016
$ws.Range("T5").Value = 0.0004837161793899014
$ws.Range("G6").Value = 1.424164
$ws.Range("H6").Value = 4.272492
$ws.Range("I6").Value = 0.03823384556371837
$ws.Range("J6").Value = 0.03823384556371837
$ws.Range("M6").Value = 87.516249
$ws.Range("N6").Value = 262.548747
$ws.Range("O6").Value = 0.3906719906954657
$ws.Range("P6").Value = 0.3906719906954657
$ws.Range("Q6").Value = 124.637491240836
$ws.Range("R6").Value = 1121.737421167524
$ws.Range("S6").Value = 0.01493689255832086
$ws.Range("T6").Value = 0.01493689255832085
$ws.Range("I7").Value = 0.08783524098133262
$ws.Range("J7").Value = 0.08783524098133261
$ws.Range("M7").Value = 19.42991633333333
$ws.Range("N7").Value = 58.289749
$ws.Range("O7").Value = 0.08673502554925175
$ws.Range("P7").Value = 0.08673502554925173
$ws.Range("Q7").Value = 63.56993238980377
$ws.Range("R7").Value = 572.129391508234
$ws.Range("S7").Value = 0.007618391870640569
$ws.Range("T7").Value = 0.007618391870640567
$ws.Range("I8").Value = 0.08783524098133262
$ws.Range("J8").Value = 0.08783524098133261
$ws.Range("O8").Value = 0.04103322570207864
$ws.Range("P8").Value = 0.04103322570207864
$ws.Range("S8").Value = 0.003604163267783489
$ws.Range("T8").Value = 0.003604163267783488
$ws.Range("I9").Value = 0.08783524098133262
$ws.Range("J9").Value = 0.08783524098133261
$ws.Range("M9").Value = 105.042315
$ws.Range("N9").Value = 315.126945
$ws.Range("O9").Value = 0.4689082402093144
$ws.Range("P9").Value = 0.4689082402093144
$ws.Range("Q9").Value = 343.6727543269299
$ws.Range("R9").Value = 3093.05478894237
$ws.Range("S9").Value = 0.04118666827691773
$ws.Range("T9").Value = 0.04118666827691773
$ws.Range("I10").Value = 0.08783524098133262
$ws.Range("J10").Value = 0.08783524098133261
$ws.Range("M10").Value = 2.834125333333334
$ws.Range("N10").Value = 8.502376000000002
$ws.Range("O10").Value = 0.01265151784388958
$ws.Range("P10").Value = 0.01265151784388958
$ws.Range("Q10").Value = 9.272564674668446
$ws.Range("R10").Value = 83.45308207201602
$ws.Range("S10").Value = 0.001111249118597671
$ws.Range("T10").Value = 0.001111249118597671
$ws.Range("I11").Value = 0.08783524098133262
$ws.Range("J11").Value = 0.08783524098133261
$ws.Range("M11").Value = 87.516249
$ws.Range("N11").Value = 262.548747
$ws.Range("O11").Value = 0.3906719906954657
$ws.Range("P11").Value = 0.3906719906954657
$ws.Range("Q11").Value = 286.331754419078
$ws.Range("R11").Value = 2576.985789771702
$ws.Range("S11").Value = 0.03431476844739317
$ws.Range("T11").Value = 0.03431476844739315
$ws.Range("G12").Value = 17.63507366666667
$ws.Range("H12").Value = 52.905221
$ws.Range("I12").Value = 0.4734403362787783
$ws.Range("J12").Value = 0.4734403362787782
$ws.Range("M12").Value = 19.42991633333333
$ws.Range("N12").Value = 58.289749
$ws.Range("O12").Value = 0.08673502554925175
$ws.Range("P12").Value = 0.08673502554925173
$ws.Range("Q12").Value = 342.6480058755033
$ws.Range("R12").Value = 3083.832052879529
$ws.Range("S12").Value = 0.04106385966318617
$ws.Range("T12").Value = 0.04106385966318617
$ws.Range("G13").Value = 17.63507366666667
$ws.Range("H13").Value = 52.905221
$ws.Range("I13").Value = 0.4734403362787783
$ws.Range("J13").Value = 0.4734403362787782
$ws.Range("O13").Value = 0.04103322570207864
$ws.Range("P13").Value = 0.04103322570207864
$ws.Range("Q13").Value = 162.1023672088834
$ws.Range("R13").Value = 1458.921304879951
$ws.Range("S13").Value = 0.01942678417499512
$ws.Range("T13").Value = 0.01942678417499512
$ws.Range("G14").Value = 17.63507366666667
$ws.Range("H14").Value = 52.905221
$ws.Range("I14").Value = 0.4734403362787783
$ws.Range("J14").Value = 0.4734403362787782
$ws.Range("M14").Value = 105.042315
$ws.Range("N14").Value = 315.126945
$ws.Range("O14").Value = 0.4689082402093144
$ws.Range("P14").Value = 0.4689082402093144
$ws.Range("Q14").Value = 1852.428963142205
$ws.Range("R14").Value = 16671.86066827984
$ws.Range("S14").Value = 0.222000074928588
$ws.Range("T14").Value = 0.2220000749285879
$ws.Range("G15").Value = 17.63507366666667
$ws.Range("H15").Value = 52.905221
$ws.Range("I15").Value = 0.4734403362787783
$ws.Range("J15").Value = 0.4734403362787782
$ws.Range("M15").Value = 2.834125333333334
$ws.Range("N15").Value = 8.502376000000002
$ws.Range("O15").Value = 0.01265151784388958
$ws.Range("P15").Value = 0.01265151784388958
$ws.Range("Q15").Value = 49.98000903389957
$ws.Range("R15").Value = 449.8200813050961
$ws.Range("S15").Value = 0.005989738862448049
$ws.Range("T15").Value = 0.005989738862448047
$ws.Range("G16").Value = 17.63507366666667
$ws.Range("H16").Value = 52.905221
$ws.Range("I16").Value = 0.4734403362787783
$ws.Range("J16").Value = 0.4734403362787782
$ws.Range("M16").Value = 87.516249
$ws.Range("N16").Value = 262.548747
$ws.Range("O16").Value = 0.3906719906954657
$ws.Range("P16").Value = 0.3906719906954657
$ws.Range("Q16").Value = 1543.355498145343
$ws.Range("R16").Value = 13890.19948330809
$ws.Range("S16").Value = 0.184959878649561
$ws.Range("T16").Value = 0.184959878649561
$ws.Range("G17").Value = 0.5460243333333333
$ws.Range("H17").Value = 1.638073
$ws.Range("I17").Value = 0.01465885251607185
$ws.Range("J17").Value = 0.01465885251607185
$ws.Range("M17").Value = 19.42991633333333
$ws.Range("N17").Value = 58.289749
$ws.Range("O17").Value = 0.08673502554925175
$ws.Range("P17").Value = 0.08673502554925173
$ws.Range("Q17").Value = 10.60920711263078
$ws.Range("R17").Value = 95.482864013677
$ws.Range("S17").Value = 0.001271435947504205
$ws.Range("T17").Value = 0.001271435947504205
$ws.Range("G18").Value = 0.5460243333333333
$ws.Range("H18").Value = 1.638073
$ws.Range("I18").Value = 0.01465885251607185
$ws.Range("J18").Value = 0.01465885251607185
$ws.Range("O18").Value = 0.04103322570207864
$ws.Range("P18").Value = 0.04103322570207864
$ws.Range("Q18").Value = 5.019079515062556
$ws.Range("R18").Value = 45.171715635563
$ws.Range("S18").Value = 0.0006015000038254596
$ws.Range("T18").Value = 0.0006015000038254595
$ws.Range("G19").Value = 0.5460243333333333
$ws.Range("H19").Value = 1.638073
$ws.Range("I19").Value = 0.01465885251607185
$ws.Range("J19").Value = 0.01465885251607185
$ws.Range("M19").Value = 105.042315
$ws.Range("N19").Value = 315.126945
$ws.Range("O19").Value = 0.4689082402093144
$ws.Range("P19").Value = 0.4689082402093144
$ws.Range("Q19").Value = 57.355660019665
$ws.Range("R19").Value = 516.200940176985
$ws.Range("S19").Value = 0.006873656736799132
$ws.Range("T19").Value = 0.006873656736799131
$ws.Range("G20").Value = 0.5460243333333333
$ws.Range("H20").Value = 1.638073
$ws.Range("I20").Value = 0.01465885251607185
$ws.Range("J20").Value = 0.01465885251607185
$ws.Range("M20").Value = 2.834125333333334
$ws.Range("N20").Value = 8.502376000000002
$ws.Range("O20").Value = 0.01265151784388958
$ws.Range("P20").Value = 0.01265151784388958
$ws.Range("Q20").Value = 1.547501395716445
$ws.Range("R20").Value = 13.927512561448
$ws.Range("S20").Value = 0.0001854567341780287
$ws.Range("T20").Value = 0.0001854567341780287
$ws.Range("G21").Value = 0.5460243333333333
$ws.Range("H21").Value = 1.638073
$ws.Range("I21").Value = 0.01465885251607185
$ws.Range("J21").Value = 0.01465885251607185
$ws.Range("M21").Value = 87.516249
$ws.Range("N21").Value = 262.548747
$ws.Range("O21").Value = 0.3906719906954657
$ws.Range("P21").Value = 0.3906719906954657
$ws.Range("Q21").Value = 47.786001516059
$ws.Range("R21").Value = 430.074013644531
$ws.Range("S21").Value = 0.005726803093765025
$ws.Range("T21").Value = 0.005726803093765024
$ws.Range("G22").Value = 14.37176
$ws.Range("H22").Value = 43.11528
$ws.Range("I22").Value = 0.385831724660099
$ws.Range("J22").Value = 0.3858317246600989
$ws.Range("M22").Value = 19.42991633333333
$ws.Range("N22").Value = 58.289749
$ws.Range("O22").Value = 0.08673502554925175
$ws.Range("P22").Value = 0.08673502554925173
$ws.Range("Q22").Value = 279.2420943627467
$ws.Range("R22").Value = 2513.17884926472
$ws.Range("S22").Value = 0.03346512449610555
$ws.Range("T22").Value = 0.03346512449610554
$ws.Range("G23").Value = 14.37176
$ws.Range("H23").Value = 43.11528
$ws.Range("I23").Value = 0.385831724660099
$ws.Range("J23").Value = 0.3858317246600989
$ws.Range("O23").Value = 0.04103322570207864
$ws.Range("P23").Value = 0.04103322570207864
$ws.Range("Q23").Value = 132.1058454868533
$ws.Range("R23").Value = 1188.95260938168
$ws.Range("S23").Value = 0.0158319202410001
$ws.Range("T23").Value = 0.0158319202410001
$ws.Range("G24").Value = 14.37176
$ws.Range("H24").Value = 43.11528
$ws.Range("I24").Value = 0.385831724660099
$ws.Range("J24").Value = 0.3858317246600989
$ws.Range("M24").Value = 105.042315
$ws.Range("N24").Value = 315.126945
$ws.Range("O24").Value = 0.4689082402093144
$ws.Range("P24").Value = 0.4689082402093144
$ws.Range("Q24").Value = 1509.6429410244
$ws.Range("R24").Value = 13586.7864692196
$ws.Range("S24").Value = 0.1809196750272918
$ws.Range("T24").Value = 0.1809196750272917
$ws.Range("G25").Value = 14.37176
$ws.Range("H25").Value = 43.11528
$ws.Range("I25").Value = 0.385831724660099
$ws.Range("J25").Value = 0.3858317246600989
$ws.Range("M25").Value = 2.834125333333334
$ws.Range("N25").Value = 8.502376000000002
$ws.Range("O25").Value = 0.01265151784388958
$ws.Range("P25").Value = 0.01265151784388958
$ws.Range("Q25").Value = 40.73136910058668
$ws.Range("R25").Value = 366.5823219052801
$ws.Range("S25").Value = 0.004881356949275935
$ws.Range("T25").Value = 0.004881356949275933
$ws.Range("G26").Value = 14.37176
$ws.Range("H26").Value = 43.11528
$ws.Range("I26").Value = 0.385831724660099
$ws.Range("J26").Value = 0.3858317246600989
$ws.Range("M26").Value = 87.516249
$ws.Range("N26").Value = 262.548747
$ws.Range("O26").Value = 0.3906719906954657
$ws.Range("P26").Value = 0.3906719906954657
$ws.Range("Q26").Value = 1257.76252672824
$ws.Range("R26").Value = 11319.86274055416
$ws.Range("S26").Value = 0.1507336479464257
$ws.Range("T26").Value = 0.1507336479464256
